$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "TestCaseName"
$ws.Range("B1").Value = "RunMode"
$ws.Range("C1").Value = "Browser"
$ws.Range("D1").Value = "AppURL"
$ws.Range("E1").Value = "UserName"
$ws.Range("F1").Value = "Password"

# --- Data rows ---
$ws.Range("A2").Value = "TC3_twitLoginChrome"
$ws.Range("B2").Value = "Y"
$ws.Range("C2").Value = "Chrome"
$ws.Range("D2").Value = "https://twitter.com/login?lang=en"
$ws.Range("E2").Value = "s1"
$ws.Range("F2").Value = "p1"

$ws.Range("A3").Value = "TC3_twitLoginChrome"
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "Chrome"
$ws.Range("D3").Value = "https://twitter.com/login?lang=en"
$ws.Range("E3").Value = "s2"
$ws.Range("F3").Value = "p2"

$ws.Range("A4").Value = "DatadrivenTest"
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "IE"
$ws.Range("D4").Value = "https://twitter.com/login?lang=en"
$ws.Range("E4").Value = "s3"
$ws.Range("F4").Value = "p3"

$ws.Range("A5").Value = "DatadrivenTest"
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Value = "Chrome"
$ws.Range("D5").Value = "https://twitter.com/login?lang=en"
$ws.Range("E5").Value = "s3"
$ws.Range("F5").Value = "p3"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.7109375
$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 8.28515625
$ws.Columns.Item(4).ColumnWidth = 32.28515625
$ws.Columns.Item(5).ColumnWidth = 10.42578125
$ws.Columns.Item(6).ColumnWidth = 9.42578125
